$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.437.41"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "2.402.40"
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.60"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.61"
$ws.Range("E6").Value = "  -2.77%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -10.47%  "
$ws.Range("D9").Value = "2.400.88"
$ws.Range("E9").Value = "  -3.26%  "
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.63"
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("D15").Value = "2.831.86"
$ws.Range("E15").Value = "  -3.36%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "60.742.05"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "2.398.07"
$ws.Range("E18").Value = "  -2.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.85"
$ws.Range("E19").Value = "  -2.83%  "
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.13"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.76"
$ws.Range("E22").Value = "  -3.64%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +6.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.87"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.28"
$ws.Range("E26").Value = "  +8.56%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0947"
$ws.Range("E28").Value = "  -4.18%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.516.88"
$ws.Range("E29").Value = "  -3.48%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "537.57"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -4.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.16"
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("E33").Value = "  -3.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("E38").Value = "  -5.14%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.379"
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.86"
$ws.Range("E40").Value = "  +7.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.20"
$ws.Range("E41").Value = "  -1.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "138.42"
$ws.Range("E42").Value = "  -7.14%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.31"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  -5.72%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.65"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.03"
$ws.Range("E47").Value = "  -4.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.45"
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0524"
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.581"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("E51").Value = "  -0.07%  "
